# Actualización automática 2025-08-04 17:26:10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Adjust column widths (D, E, F)
# Note: this runtime's saved XML "width" = ColumnWidth + 5/6 (character-width
# padding round-trip), so back the target width off by 5/6 before assigning.
$ws.Columns.Item(4).ColumnWidth = 11 - 5/6
$ws.Columns.Item(5).ColumnWidth = 22 - 5/6
$ws.Columns.Item(6).ColumnWidth = 18 - 5/6

# Row 3 - 240X80 PORCELANATO
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 4168.07156573679
$ws.Cells.Item(3, 6).Value = 0

# Row 4 - FREGADEROS DE COCINA
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 513.831046659336
$ws.Cells.Item(4, 6).Value = 0

# Row 7 - INODOROS
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 2400
$ws.Cells.Item(7, 6).Value = 0

# Row 8 - LAVABOS
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 625
$ws.Cells.Item(8, 6).Value = 0

# Row 10 - NO RESURTIBLES
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 650.25
$ws.Cells.Item(10, 6).Value = 0

# Row 12 - PANELES DECORATIVOS
$ws.Cells.Item(12, 3).Value = 100
$ws.Cells.Item(12, 5).Value = 100

# Row 13 - PANELES PU
$ws.Cells.Item(13, 3).Value = 20
$ws.Cells.Item(13, 5).Value = 20

# Row 14 - PANELES PVC
$ws.Cells.Item(14, 3).Value = 100
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 100
$ws.Cells.Item(14, 6).Value = 0

# Row 15 - PIEDRA SINTERIZADA
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 7465
$ws.Cells.Item(15, 6).Value = 0

# Row 16 - PORCELANATO
$ws.Cells.Item(16, 3).Value = 38776.47
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 38776.47
$ws.Cells.Item(16, 6).Value = 0

# Row 18 - SAL SOLUBLE
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 2800
$ws.Cells.Item(18, 6).Value = 0

# Row 19 - TOTAL
$ws.Cells.Item(19, 3).Value = 59388.22762291769
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 59388.22762291769
$ws.Cells.Item(19, 6).Value = 0
